# Update the SAS token used in the azcopy command line, and move the
# "_GoBack" bookmark (Word's "last edit location" marker) to sit right
# after the newly-typed text, the way Word itself would after a manual
# in-place edit.

$d = $word.ActiveDocument

# The old SAS token is split across two runs ("...sig=jqYQlGz...kiL" and
# "hEoo6k0...E0U%3D"), immediately followed by a run containing a single
# trailing space. Matching all of that (including the trailing space) and
# replacing it with the new token (no trailing space) reproduces a normal
# "select old text, type new text" edit.
$oldToken = "2016-06-05T21%3A27%3A00Z&se=2016-06-06T21%3A27%3A00Z&sp=rl&sv=2015-04-05&sr=c&sig=jqYQlGz%2F9fuoXDA6KiLhEoo6k0%2BD6aoRd5gH%2FkQ8E0U%3D "
$newToken = "2016-06-08T14%3A04%3A00Z&se=2016-07-06T05%3A00%3A00Z&sp=rl&sv=2015-04-05&sr=c&sig=wv9KgQfxgvkfczb655XnVLmlTu%2B5DPkKziwWC3ViWj0%3D"

$rng = $d.Content
$found = $rng.Find.Execute($oldToken, $true, $false, $false, $false, $false, $true, 1, $false, $newToken, 2)

$startPos = $rng.Start
$endPos = $rng.End

# The replace above naturally merges with whatever text/runs sit right
# before and after it. Re-establish the original run boundaries on both
# sides (the "azcopy ... ?st=" run before, and the closing quote / "/S /Y"
# run after) by briefly dropping a bookmark at each boundary -- adding and
# immediately deleting a bookmark splits the run at that point without
# leaving any trace behind.
$before = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_TmpSplitBefore", $before)
$d.Bookmarks("_TmpSplitBefore").Delete()

$after = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_TmpSplitAfter", $after)
$d.Bookmarks("_TmpSplitAfter").Delete()

# Move "_GoBack" (Word auto-maintains this bookmark at the site of the
# most recent edit) from its old location near the end of the document to
# right after the text we just typed.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$goBackRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
